# Insert a leading space before the inline picture in the first paragraph.
#
# The first paragraph currently contains only the inline picture (run with
# <w:drawing>). We need a new run containing a single space " " placed
# immediately before the drawing run, inside the same paragraph.
#
# Directly inserting text at the very start of a paragraph that begins with
# an inline shape can clobber the shape's run, so instead we:
#   1. Split the paragraph by inserting a new paragraph break right before
#      the picture (this creates an empty paragraph ahead of the picture
#      paragraph, leaving the picture's run completely untouched).
#   2. Type the space into that new (now first) paragraph.
#   3. Delete the paragraph mark that separates the space paragraph from the
#      picture paragraph, merging them back into a single paragraph that
#      reads: [space run][drawing run].

$d = $word.ActiveDocument

$firstPara = $d.Paragraphs.Item(1)
$splitPoint = $firstPara.Range.Duplicate
$splitPoint.Collapse(1)
$splitPoint.InsertParagraphBefore()

$spaceParaRange = $d.Paragraphs.Item(1).Range
$spaceParaRange.InsertBefore(" ")

$mergeStart = $spaceParaRange.Start + 1
$mergeRange = $d.Range($mergeStart, $mergeStart + 1)
$mergeRange.Delete()
